$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = "BTC"
$ws.Cells.Item(2, 3).Value = "Bitcoin"
$ws.Cells.Item(2, 4).Value = 29325
$ws.Cells.Item(2, 5).Value = 570032056521
$ws.Cells.Item(2, 6).Value = 5242729362
$ws.Cells.Item(2, 7).Value = 0.1717

$ws.Cells.Item(3, 2).Value = "ETH"
$ws.Cells.Item(3, 3).Value = "Ethereum"
$ws.Cells.Item(3, 4).Value = 1876.16
$ws.Cells.Item(3, 5).Value = 225431593307
$ws.Cells.Item(3, 6).Value = 3230890677
$ws.Cells.Item(3, 7).Value = 0.25618

$ws.Cells.Item(4, 2).Value = "USDT"
$ws.Cells.Item(4, 3).Value = "Tether"
$ws.Cells.Item(4, 4).Value = 0.999892
$ws.Cells.Item(4, 5).Value = 83816533176
$ws.Cells.Item(4, 6).Value = 7584568985
$ws.Cells.Item(4, 7).Value = 0.00054

$ws.Cells.Item(5, 2).Value = "XRP"
$ws.Cells.Item(5, 3).Value = "XRP"
$ws.Cells.Item(5, 4).Value = 0.720912
$ws.Cells.Item(5, 5).Value = 37947659364
$ws.Cells.Item(5, 6).Value = 787121276
$ws.Cells.Item(5, 7).Value = 1.62797

$ws.Cells.Item(6, 2).Value = "BNB"
$ws.Cells.Item(6, 3).Value = "BNB"
$ws.Cells.Item(6, 4).Value = 242.31
$ws.Cells.Item(6, 5).Value = 37274007780
$ws.Cells.Item(6, 6).Value = 326225353
$ws.Cells.Item(6, 7).Value = 0.29641

$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "USD Coin"
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 26552861517
$ws.Cells.Item(7, 6).Value = 1678874930
$ws.Cells.Item(7, 7).Value = -0.02898

$ws.Cells.Item(8, 2).Value = "STETH"
$ws.Cells.Item(8, 3).Value = "Lido Staked Ether"
$ws.Cells.Item(8, 4).Value = 1875.71
$ws.Cells.Item(8, 5).Value = 14861884752
$ws.Cells.Item(8, 6).Value = 15262107
$ws.Cells.Item(8, 7).Value = 0.28481

$ws.Cells.Item(9, 2).Value = "DOGE"
$ws.Cells.Item(9, 3).Value = "Dogecoin"
$ws.Cells.Item(9, 4).Value = 0.080179
$ws.Cells.Item(9, 5).Value = 11248371166
$ws.Cells.Item(9, 6).Value = 601097859
$ws.Cells.Item(9, 7).Value = 2.6147

$ws.Cells.Item(10, 2).Value = "ADA"
$ws.Cells.Item(10, 3).Value = "Cardano"
$ws.Cells.Item(10, 4).Value = 0.31423
$ws.Cells.Item(10, 5).Value = 11014127034
$ws.Cells.Item(10, 6).Value = 144511044
$ws.Cells.Item(10, 7).Value = 1.52226

$ws.Cells.Item(11, 2).Value = "SOL"
$ws.Cells.Item(11, 3).Value = "Solana"
$ws.Cells.Item(11, 4).Value = 24.98
$ws.Cells.Item(11, 5).Value = 10093471873
$ws.Cells.Item(11, 6).Value = 245013141
$ws.Cells.Item(11, 7).Value = -0.14491

$ws.Cells.Item(12, 2).Value = "TRX"
$ws.Cells.Item(12, 3).Value = "TRON"
$ws.Cells.Item(12, 4).Value = 0.081881
$ws.Cells.Item(12, 5).Value = 7342197506
$ws.Cells.Item(12, 6).Value = 233944871
$ws.Cells.Item(12, 7).Value = -2.5746

$ws.Cells.Item(13, 2).Value = "LTC"
$ws.Cells.Item(13, 3).Value = "Litecoin"
$ws.Cells.Item(13, 4).Value = 94.56
$ws.Cells.Item(13, 5).Value = 6945205662
$ws.Cells.Item(13, 6).Value = 971326813
$ws.Cells.Item(13, 7).Value = 3.90549

$ws.Cells.Item(14, 2).Value = "MATIC"
$ws.Cells.Item(14, 3).Value = "Polygon"
$ws.Cells.Item(14, 4).Value = 0.712296
$ws.Cells.Item(14, 5).Value = 6635121757
$ws.Cells.Item(14, 6).Value = 133360835
$ws.Cells.Item(14, 7).Value = 0.2121

$ws.Cells.Item(15, 2).Value = "DOT"
$ws.Cells.Item(15, 3).Value = "Polkadot"
$ws.Cells.Item(15, 4).Value = 5.23
$ws.Cells.Item(15, 5).Value = 6577578181
$ws.Cells.Item(15, 6).Value = 72966184
$ws.Cells.Item(15, 7).Value = -0.09211

$ws.Cells.Item(16, 2).Value = "SHIB"
$ws.Cells.Item(16, 3).Value = "Shiba Inu"
$ws.Cells.Item(16, 4).Value = 0.00000851
$ws.Cells.Item(16, 5).Value = 5018654493
$ws.Cells.Item(16, 6).Value = 243230312
$ws.Cells.Item(16, 7).Value = 3.90805

$ws.Cells.Item(17, 2).Value = "UNI"
$ws.Cells.Item(17, 3).Value = "Uniswap"
$ws.Cells.Item(17, 4).Value = 6.4
$ws.Cells.Item(17, 5).Value = 4826523067
$ws.Cells.Item(17, 6).Value = 148256215
$ws.Cells.Item(17, 7).Value = 5.48412

$ws.Cells.Item(18, 2).Value = "WBTC"
$ws.Cells.Item(18, 3).Value = "Wrapped Bitcoin"
$ws.Cells.Item(18, 4).Value = 29357
$ws.Cells.Item(18, 5).Value = 4738917991
$ws.Cells.Item(18, 6).Value = 24063273
$ws.Cells.Item(18, 7).Value = 0.13688

$ws.Cells.Item(19, 2).Value = "BCH"
$ws.Cells.Item(19, 3).Value = "Bitcoin Cash"
$ws.Cells.Item(19, 4).Value = 243.48
$ws.Cells.Item(19, 5).Value = 4730881411
$ws.Cells.Item(19, 6).Value = 117413214
$ws.Cells.Item(19, 7).Value = 1.507

$ws.Cells.Item(20, 2).Value = "AVAX"
$ws.Cells.Item(20, 3).Value = "Avalanche"
$ws.Cells.Item(20, 4).Value = 13.29
$ws.Cells.Item(20, 5).Value = 4595772815
$ws.Cells.Item(20, 6).Value = 64229971
$ws.Cells.Item(20, 7).Value = 0.63089

$ws.Cells.Item(21, 2).Value = "XLM"
$ws.Cells.Item(21, 3).Value = "Stellar"
$ws.Cells.Item(21, 4).Value = 0.1596
$ws.Cells.Item(21, 5).Value = 4351163512
$ws.Cells.Item(21, 6).Value = 121616276
$ws.Cells.Item(21, 7).Value = 0.10115

$ws.Cells.Item(22, 2).Value = "LINK"
$ws.Cells.Item(22, 3).Value = "Chainlink"
$ws.Cells.Item(22, 4).Value = 7.75
$ws.Cells.Item(22, 5).Value = 4169407411
$ws.Cells.Item(22, 6).Value = 131714157
$ws.Cells.Item(22, 7).Value = 0.14586

$ws.Cells.Item(23, 2).Value = "TON"
$ws.Cells.Item(23, 3).Value = "Toncoin"
$ws.Cells.Item(23, 4).Value = 1.21
$ws.Cells.Item(23, 5).Value = 4162500783
$ws.Cells.Item(23, 6).Value = 33747336
$ws.Cells.Item(23, 7).Value = -6.56901

$ws.Cells.Item(24, 2).Value = "DAI"
$ws.Cells.Item(24, 3).Value = "Dai"
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 4151401398
$ws.Cells.Item(24, 6).Value = 46218550
$ws.Cells.Item(24, 7).Value = 0.01485

$ws.Cells.Item(25, 2).Value = "LEO"
$ws.Cells.Item(25, 3).Value = "LEO Token"
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = 3725423698
$ws.Cells.Item(25, 6).Value = 1043572
$ws.Cells.Item(25, 7).Value = 1.52728

$ws.Cells.Item(26, 2).Value = "BUSD"
$ws.Cells.Item(26, 3).Value = "Binance USD"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 3705611563
$ws.Cells.Item(26, 6).Value = 712448899
$ws.Cells.Item(26, 7).Value = 0.0121

$ws.Cells.Item(27, 2).Value = "TUSD"
$ws.Cells.Item(27, 3).Value = "TrueUSD"
$ws.Cells.Item(27, 4).Value = 0.998335
$ws.Cells.Item(27, 5).Value = 3053818573
$ws.Cells.Item(27, 6).Value = 74033362
$ws.Cells.Item(27, 7).Value = -0.09396

$ws.Cells.Item(28, 2).Value = "XMR"
$ws.Cells.Item(28, 3).Value = "Monero"
$ws.Cells.Item(28, 4).Value = 162.46
$ws.Cells.Item(28, 5).Value = 2944328612
$ws.Cells.Item(28, 6).Value = 54333416
$ws.Cells.Item(28, 7).Value = -0.09002

$ws.Cells.Item(29, 2).Value = "ATOM"
$ws.Cells.Item(29, 3).Value = "Cosmos Hub"
$ws.Cells.Item(29, 4).Value = 9.029999999999999
$ws.Cells.Item(29, 5).Value = 2641700611
$ws.Cells.Item(29, 6).Value = 54545411
$ws.Cells.Item(29, 7).Value = 0.27636

$ws.Cells.Item(30, 2).Value = "ETC"
$ws.Cells.Item(30, 3).Value = "Ethereum Classic"
$ws.Cells.Item(30, 4).Value = 18.48
$ws.Cells.Item(30, 5).Value = 2629216357
$ws.Cells.Item(30, 6).Value = 57201279
$ws.Cells.Item(30, 7).Value = 0.24892

$ws.Cells.Item(31, 2).Value = "OKB"
$ws.Cells.Item(31, 3).Value = "OKB"
$ws.Cells.Item(31, 4).Value = 42.73
$ws.Cells.Item(31, 5).Value = 2564141781
$ws.Cells.Item(31, 6).Value = 1395193
$ws.Cells.Item(31, 7).Value = 0.14302

$ws.Cells.Item(32, 2).Value = "FIL"
$ws.Cells.Item(32, 3).Value = "Filecoin"
$ws.Cells.Item(32, 4).Value = 4.4
$ws.Cells.Item(32, 5).Value = 1929573013
$ws.Cells.Item(32, 6).Value = 52314459
$ws.Cells.Item(32, 7).Value = 0.30659

$ws.Cells.Item(33, 2).Value = "ICP"
$ws.Cells.Item(33, 3).Value = "Internet Computer"
$ws.Cells.Item(33, 4).Value = 4.29
$ws.Cells.Item(33, 5).Value = 1887796438
$ws.Cells.Item(33, 6).Value = 12690369
$ws.Cells.Item(33, 7).Value = -0.1649

$ws.Cells.Item(34, 2).Value = "HBAR"
$ws.Cells.Item(34, 3).Value = "Hedera"
$ws.Cells.Item(34, 4).Value = 0.053558
$ws.Cells.Item(34, 5).Value = 1735835648
$ws.Cells.Item(34, 6).Value = 19061521
$ws.Cells.Item(34, 7).Value = 0.04295

$ws.Cells.Item(35, 2).Value = "LDO"
$ws.Cells.Item(35, 3).Value = "Lido DAO"
$ws.Cells.Item(35, 4).Value = 1.93
$ws.Cells.Item(35, 5).Value = 1698206638
$ws.Cells.Item(35, 6).Value = 20429035
$ws.Cells.Item(35, 7).Value = -0.47739

$ws.Cells.Item(36, 2).Value = "MNT"
$ws.Cells.Item(36, 3).Value = "Mantle"
$ws.Cells.Item(36, 4).Value = 0.521925
$ws.Cells.Item(36, 5).Value = 1687718650
$ws.Cells.Item(36, 6).Value = 6546496
$ws.Cells.Item(36, 7).Value = 0.4054

$ws.Cells.Item(37, 2).Value = "QNT"
$ws.Cells.Item(37, 3).Value = "Quant"
$ws.Cells.Item(37, 4).Value = 112.88
$ws.Cells.Item(37, 5).Value = 1640994986
$ws.Cells.Item(37, 6).Value = 19975585
$ws.Cells.Item(37, 7).Value = 3.84517

$ws.Cells.Item(38, 2).Value = "CRO"
$ws.Cells.Item(38, 3).Value = "Cronos"
$ws.Cells.Item(38, 4).Value = 0.059765
$ws.Cells.Item(38, 5).Value = 1570026240
$ws.Cells.Item(38, 6).Value = 3760921
$ws.Cells.Item(38, 7).Value = -0.24554

$ws.Cells.Item(39, 2).Value = "APT"
$ws.Cells.Item(39, 3).Value = "Aptos"
$ws.Cells.Item(39, 4).Value = 7.09
$ws.Cells.Item(39, 5).Value = 1546432165
$ws.Cells.Item(39, 6).Value = 26643714
$ws.Cells.Item(39, 7).Value = 0.34894

$ws.Cells.Item(40, 2).Value = "ARB"
$ws.Cells.Item(40, 3).Value = "Arbitrum"
$ws.Cells.Item(40, 4).Value = 1.18
$ws.Cells.Item(40, 5).Value = 1500095010
$ws.Cells.Item(40, 6).Value = 53154525
$ws.Cells.Item(40, 7).Value = 0.12822

$ws.Cells.Item(41, 2).Value = "VET"
$ws.Cells.Item(41, 3).Value = "VeChain"
$ws.Cells.Item(41, 4).Value = 0.01874645
$ws.Cells.Item(41, 5).Value = 1362328054
$ws.Cells.Item(41, 6).Value = 27407729
$ws.Cells.Item(41, 7).Value = 0.28904

$ws.Cells.Item(42, 2).Value = "NEAR"
$ws.Cells.Item(42, 3).Value = "NEAR Protocol"
$ws.Cells.Item(42, 4).Value = 1.41
$ws.Cells.Item(42, 5).Value = 1330956238
$ws.Cells.Item(42, 6).Value = 33739112
$ws.Cells.Item(42, 7).Value = 0.8604000000000001

$ws.Cells.Item(43, 2).Value = "MKR"
$ws.Cells.Item(43, 3).Value = "Maker"
$ws.Cells.Item(43, 4).Value = 1276.84
$ws.Cells.Item(43, 5).Value = 1150872039
$ws.Cells.Item(43, 6).Value = 115632210
$ws.Cells.Item(43, 7).Value = 4.74207

$ws.Cells.Item(44, 2).Value = "OP"
$ws.Cells.Item(44, 3).Value = "Optimism"
$ws.Cells.Item(44, 4).Value = 1.58
$ws.Cells.Item(44, 5).Value = 1074440841
$ws.Cells.Item(44, 6).Value = 98374366
$ws.Cells.Item(44, 7).Value = 5.86569

$ws.Cells.Item(45, 2).Value = "AAVE"
$ws.Cells.Item(45, 3).Value = "Aave"
$ws.Cells.Item(45, 4).Value = 74.15000000000001
$ws.Cells.Item(45, 5).Value = 1073493413
$ws.Cells.Item(45, 6).Value = 49869146
$ws.Cells.Item(45, 7).Value = 2.62686

$ws.Cells.Item(46, 2).Value = "GRT"
$ws.Cells.Item(46, 3).Value = "The Graph"
$ws.Cells.Item(46, 4).Value = 0.11188
$ws.Cells.Item(46, 5).Value = 1019200908
$ws.Cells.Item(46, 6).Value = 11864905
$ws.Cells.Item(46, 7).Value = -0.03139

$ws.Cells.Item(47, 2).Value = "RETH"
$ws.Cells.Item(47, 3).Value = "Rocket Pool ETH"
$ws.Cells.Item(47, 4).Value = 2022.88
$ws.Cells.Item(47, 5).Value = 985444103
$ws.Cells.Item(47, 6).Value = 1796287
$ws.Cells.Item(47, 7).Value = 0.28585

$ws.Cells.Item(48, 2).Value = "FRAX"
$ws.Cells.Item(48, 3).Value = "Frax"
$ws.Cells.Item(48, 4).Value = 0.9993069999999999
$ws.Cells.Item(48, 5).Value = 911805361
$ws.Cells.Item(48, 6).Value = 3186966
$ws.Cells.Item(48, 7).Value = 0.11279

$ws.Cells.Item(49, 2).Value = "SNX"
$ws.Cells.Item(49, 3).Value = "Synthetix Network"
$ws.Cells.Item(49, 4).Value = 2.81
$ws.Cells.Item(49, 5).Value = 900803317
$ws.Cells.Item(49, 6).Value = 38084042
$ws.Cells.Item(49, 7).Value = -0.00086

$ws.Cells.Item(50, 2).Value = "ALGO"
$ws.Cells.Item(50, 3).Value = "Algorand"
$ws.Cells.Item(50, 4).Value = 0.110922
$ws.Cells.Item(50, 5).Value = 865397798
$ws.Cells.Item(50, 6).Value = 25471750
$ws.Cells.Item(50, 7).Value = 0.5053299999999999

$ws.Cells.Item(51, 2).Value = "EGLD"
$ws.Cells.Item(51, 3).Value = "MultiversX"
$ws.Cells.Item(51, 4).Value = 32.66
$ws.Cells.Item(51, 5).Value = 841245870
$ws.Cells.Item(51, 6).Value = 5878902
$ws.Cells.Item(51, 7).Value = -0.0592
